$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 85.7
$ws.Range("I6").Value = 99.066666
$ws.Range("J6").Value = 45.6
$ws.Range("K6").Value = 297.199998
$ws.Range("L6").Value = 136.8
$ws.Range("M6").Value = -185.199998
$ws.Range("N6").Value = -360.8
$ws.Range("H8").Value = 103.92308
$ws.Range("I8").Value = 41.6
$ws.Range("K8").Value = 124.8
$ws.Range("M8").Value = 14.19999999999999
$ws.Range("H10").Value = 27999.2
$ws.Range("I10").Value = 20000
$ws.Range("K10").Value = 20000
$ws.Range("M10").Value = -19707
$ws.Range("H17").Value = 669.9677
$ws.Range("I17").Value = 200
$ws.Range("J17").Value = 702.37933
$ws.Range("K17").Value = 600
$ws.Range("L17").Value = 2107.13799
$ws.Range("M17").Value = -432
$ws.Range("N17").Value = -2443.13799
$ws.Range("H19").Value = 3733.75
$ws.Range("I19").Value = 3309
$ws.Range("J19").Value = 3946.125
$ws.Range("K19").Value = 3309
$ws.Range("L19").Value = 3946.125
$ws.Range("M19").Value = -3134
$ws.Range("N19").Value = -4296.125
$ws.Range("H98").Value = 2911.375
$ws.Range("I98").Value = 2748.5
$ws.Range("J98").Value = 3074.25
$ws.Range("K98").Value = 2748.5
$ws.Range("L98").Value = 3074.25
$ws.Range("M98").Value = -1250.5
$ws.Range("N98").Value = -6070.25
$ws.Range("H103").Value = 511.5
$ws.Range("I103").Value = 511.15
$ws.Range("K103").Value = 1533.45
$ws.Range("M103").Value = -947.4499999999998
$ws.Range("H122").Value = 2911.375
$ws.Range("I122").Value = 2748.5
$ws.Range("J122").Value = 3074.25
$ws.Range("K122").Value = 8245.5
$ws.Range("L122").Value = 9222.75
$ws.Range("M122").Value = -5795.5
$ws.Range("N122").Value = -14122.75
$ws.Range("H129").Value = 483.5
$ws.Range("I129").Value = 483.5
$ws.Range("K129").Value = 1450.5
$ws.Range("M129").Value = 3549.5
$ws.Range("H132").Value = 50766.24
$ws.Range("I132").Value = 53154.85
$ws.Range("K132").Value = 159464.55
$ws.Range("M132").Value = -156934.55

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2573054.5
$ws.Range("I74").Value = 1425333.5
$ws.Range("J74").Value = 5557129
$ws.Range("K74").Value = 1425333.5
$ws.Range("L74").Value = 5557129
$ws.Range("M74").Value = -1424459.5
$ws.Range("N74").Value = -5558877
$ws.Range("H77").Value = 2573054.5
$ws.Range("I77").Value = 1425333.5
$ws.Range("J77").Value = 5557129
$ws.Range("K77").Value = 7126667.5
$ws.Range("L77").Value = 27785645
$ws.Range("M77").Value = -7122299.5
$ws.Range("N77").Value = -27794381
$ws.Range("H97").Value = 663.4
$ws.Range("I97").Value = 241.75
$ws.Range("K97").Value = 241.75
$ws.Range("M97").Value = 254.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1265.65
$ws.Range("I20").Value = 1156.8667
$ws.Range("K20").Value = 1156.8667
$ws.Range("M20").Value = -909.8667
$ws.Range("H99").Value = 1935.875
$ws.Range("J99").Value = 1955.1428
$ws.Range("L99").Value = 1955.1428
$ws.Range("N99").Value = -4951.1428
$ws.Range("H105").Value = 3309.9
$ws.Range("I105").Value = 3468
$ws.Range("J105").Value = 3072.75
$ws.Range("K105").Value = 3468
$ws.Range("L105").Value = 3072.75
$ws.Range("M105").Value = -1721
$ws.Range("N105").Value = -6566.75
$ws.Range("H107").Value = 58560.637
$ws.Range("I107").Value = 30702.25
$ws.Range("J107").Value = 74479.71000000001
$ws.Range("K107").Value = 30702.25
$ws.Range("L107").Value = 74479.71000000001
$ws.Range("M107").Value = -28782.25
$ws.Range("N107").Value = -78319.71000000001
$ws.Range("H135").Value = 59993
$ws.Range("J135").Value = 59993
$ws.Range("L135").Value = 59993
$ws.Range("N135").Value = -70133

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1477.1
$ws.Range("J31").Value = 1764.3334
$ws.Range("L31").Value = 1764.3334
$ws.Range("N31").Value = -2354.3334
$ws.Range("H34").Value = 1477.1
$ws.Range("J34").Value = 1764.3334
$ws.Range("L34").Value = 1764.3334
$ws.Range("N34").Value = -2168.3334

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 582.75
$ws.Range("I8").Value = 582.75
$ws.Range("K8").Value = 1748.25
$ws.Range("M8").Value = -1609.25
$ws.Range("H23").Value = 3701.875
$ws.Range("I23").Value = 4824.4287
$ws.Range("K23").Value = 14473.2861
$ws.Range("M23").Value = -14238.2861
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("N62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("N65").Value = 0
$ws.Range("H86").Value = 300
$ws.Range("H89").Value = 300
$ws.Range("H94").Value = 11498.75
$ws.Range("I94").Value = 12997.5
$ws.Range("J94").Value = 10000
$ws.Range("K94").Value = 38992.5
$ws.Range("L94").Value = 30000
$ws.Range("M94").Value = -38316.5
$ws.Range("N94").Value = -31352

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 95765.69500000001
$ws.Range("J18").Value = 24954
$ws.Range("L18").Value = 24954
$ws.Range("N18").Value = -25540
$ws.Range("H107").Value = 875.9259
$ws.Range("I107").Value = 990.5909
$ws.Range("J107").Value = 371.4
$ws.Range("K107").Value = 990.5909
$ws.Range("L107").Value = 371.4
$ws.Range("M107").Value = 929.4091
$ws.Range("N107").Value = -4211.4
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5241.4546
$ws.Range("I22").Value = 1813.5714
$ws.Range("K22").Value = 1813.5714
$ws.Range("M22").Value = -1518.5714
$ws.Range("H27").Value = 5241.4546
$ws.Range("I27").Value = 1813.5714
$ws.Range("K27").Value = 1813.5714
$ws.Range("M27").Value = -1706.5714
$ws.Range("H43").Value = 8602.4
$ws.Range("I43").Value = 5000
$ws.Range("J43").Value = 9503
$ws.Range("K43").Value = 5000
$ws.Range("L43").Value = 9503
$ws.Range("M43").Value = -4807
$ws.Range("N43").Value = -9889
$ws.Range("H122").Value = 1998.5
$ws.Range("I122").Value = 1998.5
$ws.Range("K122").Value = 5995.5
$ws.Range("M122").Value = -3545.5
$ws.Range("H132").Value = 3432.2856
$ws.Range("I132").Value = 3100.4
$ws.Range("J132").Value = 4262
$ws.Range("K132").Value = 9301.200000000001
$ws.Range("L132").Value = 12786
$ws.Range("M132").Value = -6771.200000000001
$ws.Range("N132").Value = -17846

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 4854
$ws.Range("I113").Value = 741.3333
$ws.Range("J113").Value = 6910.3335
$ws.Range("K113").Value = 2223.9999
$ws.Range("L113").Value = 20731.0005
$ws.Range("M113").Value = -53.9998999999998
$ws.Range("N113").Value = -25071.0005
$ws.Range("H132").Value = 1812.6
$ws.Range("I132").Value = 1812.6
$ws.Range("K132").Value = 5437.799999999999
$ws.Range("M132").Value = -2907.799999999999
